$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 35.333332
$ws.Range("I11").Value = 35.333332
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 35.333332
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 104.666668
$ws.Range("H33").Value = 241.46666
$ws.Range("I33").Value = 205
$ws.Range("J33").Value = 752
$ws.Range("K33").Value = 205
$ws.Range("L33").Value = 752
$ws.Range("M33").Value = 24
$ws.Range("N33").Value = -1210
$ws.Range("H39").Value = 25.666666
$ws.Range("I39").Value = 25.666666
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 76.99999800000001
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 219.000002
$ws.Range("N39").ClearContents()
$ws.Range("H53").Value = 539.6316
$ws.Range("I53").Value = 541.7222
$ws.Range("J53").Value = 502
$ws.Range("K53").Value = 541.7222
$ws.Range("L53").Value = 502
$ws.Range("M53").Value = 95.27779999999996
$ws.Range("N53").Value = -1776
$ws.Range("H94").Value = 4652.4
$ws.Range("I94").Value = 4652.4
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 4652.4
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -4201.4
$ws.Range("H141").Value = 25397.5
$ws.Range("I141").Value = 25397.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 76192.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -71012.5
$ws.Range("N141").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7105.6
$ws.Range("I32").Value = 7105.6
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 7105.6
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -6818.6
$ws.Range("H45").Value = 2425.1
$ws.Range("I45").Value = 2610.1428
$ws.Range("J45").Value = 1993.3334
$ws.Range("K45").Value = 2610.1428
$ws.Range("L45").Value = 1993.3334
$ws.Range("M45").Value = -2233.1428
$ws.Range("N45").Value = -2747.3334
$ws.Range("H63").Value = 1953
$ws.Range("I63").Value = 1700
$ws.Range("J63").Value = 2206
$ws.Range("K63").Value = 1700
$ws.Range("L63").Value = 2206
$ws.Range("M63").Value = -1014
$ws.Range("N63").Value = -3578
$ws.Range("H66").Value = 1953
$ws.Range("I66").Value = 1700
$ws.Range("J66").Value = 2206
$ws.Range("K66").Value = 8500
$ws.Range("L66").Value = 11030
$ws.Range("M66").Value = -5068
$ws.Range("N66").Value = -17894
$ws.Range("H132").Value = 2019.4
$ws.Range("I132").Value = 2019.4
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6058.200000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3528.200000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3053.5557
$ws.Range("I20").Value = 2226
$ws.Range("J20").Value = 5950
$ws.Range("K20").Value = 2226
$ws.Range("L20").Value = 5950
$ws.Range("M20").Value = -1979
$ws.Range("N20").Value = -6444
$ws.Range("H96").Value = 200000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 200000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 200000
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -205492

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H99").Value = 5574.1816
$ws.Range("I99").Value = 5078.4287
$ws.Range("J99").Value = 6441.75
$ws.Range("K99").Value = 5078.4287
$ws.Range("L99").Value = 6441.75
$ws.Range("M99").Value = -3580.4287
$ws.Range("N99").Value = -9437.75
$ws.Range("H126").Value = 5574.1816
$ws.Range("I126").Value = 5078.4287
$ws.Range("J126").Value = 6441.75
$ws.Range("K126").Value = 15235.2861
$ws.Range("L126").Value = 19325.25
$ws.Range("M126").Value = -12765.2861
$ws.Range("N126").Value = -24265.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 71.333336
$ws.Range("J2").Value = 56.75
$ws.Range("K2").Value = 428.000016
$ws.Range("L2").Value = 340.5
$ws.Range("M2").Value = -315.000016
$ws.Range("N2").Value = -566.5
$ws.Range("H11").Value = 318
$ws.Range("I11").Value = 347.5
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 1042.5
$ws.Range("L11").Value = 600
$ws.Range("M11").Value = -902.5
$ws.Range("N11").Value = -880
$ws.Range("H26").Value = 2333.3333
$ws.Range("I26").Value = 1000
$ws.Range("J26").Value = 5000
$ws.Range("K26").Value = 3000
$ws.Range("L26").Value = 15000
$ws.Range("M26").Value = -2712
$ws.Range("N26").Value = -15576
$ws.Range("H38").Value = 900.1111
$ws.Range("I38").Value = 1911
$ws.Range("J38").Value = 91.40000000000001
$ws.Range("K38").Value = 5733
$ws.Range("L38").Value = 274.2
$ws.Range("M38").Value = -5386
$ws.Range("N38").Value = -968.2
$ws.Range("H45").Value = 1500
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 4500
$ws.Range("N45").Value = -5564
$ws.Range("H68").Value = 439.75
$ws.Range("I68").Value = 442.77777
$ws.Range("J68").Value = 430.66666
$ws.Range("K68").Value = 1328.33331
$ws.Range("L68").Value = 1291.99998
$ws.Range("M68").Value = -517.33331
$ws.Range("N68").Value = -2913.99998
$ws.Range("H71").Value = 439.75
$ws.Range("I71").Value = 442.77777
$ws.Range("J71").Value = 430.66666
$ws.Range("K71").Value = 3984.99993
$ws.Range("L71").Value = 3875.99994
$ws.Range("M71").Value = 71.00007000000005
$ws.Range("N71").Value = -11987.99994
$ws.Range("H98").Value = 796.3333
$ws.Range("I98").Value = 699.5
$ws.Range("J98").Value = 990
$ws.Range("K98").Value = 2098.5
$ws.Range("L98").Value = 2970
$ws.Range("M98").Value = -600.5
$ws.Range("N98").Value = -5966
$ws.Range("H107").Value = 2357
$ws.Range("I107").Value = 2101.75
$ws.Range("J107").Value = 2561.2
$ws.Range("K107").Value = 6305.25
$ws.Range("L107").Value = 7683.599999999999
$ws.Range("M107").Value = -4385.25
$ws.Range("N107").Value = -11523.6
$ws.Range("H121").Value = 1375.15
$ws.Range("I121").Value = 1366.6666
$ws.Range("J121").Value = 1376.6471
$ws.Range("K121").Value = 4099.9998
$ws.Range("L121").Value = 4129.9413
$ws.Range("M121").Value = -2789.9998
$ws.Range("N121").Value = -6749.9413
$ws.Range("H139").Value = 4454.4546
$ws.Range("I139").Value = 4454.4546
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 13363.3638
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -8223.363799999999
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 2681.2856
$ws.Range("I140").Value = 2681.2856
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 8043.8568
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -2863.8568
$ws.Range("N140").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 22.666666
$ws.Range("I2").Value = 24
$ws.Range("J2").Value = 22
$ws.Range("K2").Value = 24
$ws.Range("L2").Value = 22
$ws.Range("M2").Value = 89
$ws.Range("N2").Value = -248
$ws.Range("H43").Value = 1651.375
$ws.Range("I43").Value = 1651.375
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1651.375
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -1500.375
$ws.Range("N43").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H95").Value = 30000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 30000
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492
$ws.Range("H132").Value = 5008
$ws.Range("I132").Value = 5008
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15024
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12494

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1326.5
$ws.Range("I82").Value = 1326.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1326.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -965.5
$ws.Range("H85").Value = 1326.5
$ws.Range("I85").Value = 1326.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1326.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -78.5
$ws.Range("H93").Value = 999.75
$ws.Range("I93").Value = 933
$ws.Range("J93").Value = 1200
$ws.Range("K93").Value = 933
$ws.Range("L93").Value = 1200
$ws.Range("M93").Value = 315
$ws.Range("N93").Value = -3696
$ws.Range("H132").Value = 2685.5
$ws.Range("I132").Value = 2685.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8056.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5526.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 100071
$ws.Range("I51").Value = 100071
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 100071
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -99561
$ws.Range("N51").ClearContents()
$ws.Range("H81").Value = 500
$ws.Range("I81").Value = 500
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 61
$ws.Range("H84").Value = 500
$ws.Range("I84").Value = 500
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 5000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 304
$ws.Range("H132").Value = 1765.1666
$ws.Range("I132").Value = 1880.1818
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 5640.5454
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -3110.5454
$ws.Range("N132").Value = -6560
